$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Extend the trailing blank (placeholder) rows.
#    The sheet used to end with blank styled rows 165-167; those
#    become real data rows 165-168, and the blank placeholder block
#    now spans rows 169-185. Clone the blank-row formatting (taken
#    from the still-untouched row 165) down across the new rows
#    before row 165 itself gets overwritten with real data.
# ------------------------------------------------------------------
$ws.Range("A165:E165").Copy()
$ws.Range("A169:E185").PasteSpecial(-4122)
$ws.Range("A169:E185").RowHeight = 15

# ------------------------------------------------------------------
# 2) Turn rows 165-168 into real data rows. Clone formatting from the
#    previous populated row (164) so the date/text cell styles match
#    the rest of the table (border/fill/number-format).
# ------------------------------------------------------------------
$ws.Range("A164:E164").Copy()
$ws.Range("A165:E168").PasteSpecial(-4122)
$ws.Range("A165:E168").RowHeight = 15

# Row 165
$ws.Cells.Item(165,1).Value = 45738
$ws.Cells.Item(165,2).Value = "HARD.NOISE HARD TECHNO"
$ws.Cells.Item(165,3).Value = "AREA 15"
$ws.Cells.Item(165,4).Value = "Bochum"
$ws.Hyperlinks.Add($ws.Cells.Item(165,5), "https://www.instagram.com/reel/DFTLQqsNP5d/?igsh=ZWRrYmhubDdtM3dy", "", "", "https://www.instagram.com/reel/DFTLQqsNP5d/?igsh=ZWRrYmhubDdtM3dy")

# Row 166
$ws.Cells.Item(166,1).Value = 45715
$ws.Cells.Item(166,2).Value = "POLAAR 180 MIN RAVE (20 Uhr)"
$ws.Cells.Item(166,3).Value = "AREA 15"
$ws.Cells.Item(166,4).Value = "Bochum"
$ws.Hyperlinks.Add($ws.Cells.Item(166,5), "https://www.instagram.com/polaartechno?igsh=MTZlYmtzODdubzVhaQ==", "", "", "https://www.instagram.com/polaartechno?igsh=MTZlYmtzODdubzVhaQ==")

# Row 167
$ws.Cells.Item(167,1).Value = 45699
$ws.Cells.Item(167,2).Value = "LOL 120 MIN RAVE (21 Uhr)"
$ws.Cells.Item(167,3).Value = "Goethebunker"
$ws.Cells.Item(167,4).Value = "Essen"
$ws.Hyperlinks.Add($ws.Cells.Item(167,5), "https://www.instagram.com/reel/DF7NXBvOScN/?igsh=MXg0cDJ0Mmk4ZTk5Mw==", "", "", "https://www.instagram.com/reel/DF7NXBvOScN/?igsh=MXg0cDJ0Mmk4ZTk5Mw==")

# Row 168
$ws.Cells.Item(168,1).Value = 45703
$ws.Cells.Item(168,2).Value = "SUPREMACY 2025 GERMANY DECODED"
$ws.Cells.Item(168,3).Value = "Westfalenhallen"
$ws.Cells.Item(168,4).Value = "Dortmund"
$ws.Hyperlinks.Add($ws.Cells.Item(168,5), "https://www.instagram.com/supremacyevent?igsh=ejB3MjBkOHRhOWxl", "", "", "https://www.instagram.com/supremacyevent?igsh=ejB3MjBkOHRhOWxl")

$ws.Range("A165:E168").RowHeight = 15
